# Edit "1 - HTML + CSS.pptx" (slide 1, body placeholder) to:
#  - turn the existing "html" run into a hyperlink (same target as the
#    "https://www.w3schools.com/" run right before it)
#  - append an "And:" line plus a second hyperlink-styled line pointing at
#    https://www.internetingishard.com/html-and-css/basic-web-pages/

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- locate the "See: <link>html" paragraph -------------------------------
$seePara = $tr.Paragraphs(3)

# 1) Insert the new paragraphs after it (done *before* the hyperlink below
#    is added, so the new runs do not inherit it):
#      (blank)
#      And:
#      (blank)
#      <url text, split into separate runs below>
$urlParaText = "https://www.internetingishard.com/html-and-css/basic-web-pages/"
$inserted = $seePara.InsertAfter("`r`rAnd:`r`r" + $urlParaText)

# Figure out where things landed.
$blank1Idx = 4
$andIdx    = 5
$blank2Idx = 6
$urlIdx    = 7

$blank1  = $tr.Paragraphs($blank1Idx)
$andPara = $tr.Paragraphs($andIdx)
$blank2  = $tr.Paragraphs($blank2Idx)
$urlPara = $tr.Paragraphs($urlIdx)

# --- formatting: blank line after "html" (red text / yellow highlight) ----
$blank1.Font.Color.RGB = 255       # FF0000
$blank1.Font.Highlight.RGB = 65535 # FFFF00

# --- formatting: "And:" (white text / black highlight) --------------------
$andPara.Font.Color.RGB = 16777215 # FFFFFF (bg1)
$andPara.Font.Highlight.RGB = 0    # 000000

# --- formatting: blank line before the link (red text / yellow highlight) -
$blank2.Font.Color.RGB = 255       # FF0000
$blank2.Font.Highlight.RGB = 65535 # FFFF00

# --- formatting: link line (dark red text / yellow highlight) -------------
$urlPara.Font.Color.RGB = 192        # C00000
$urlPara.Font.Highlight.RGB = 65535  # FFFF00

# Split the link line into the same run boundaries as the source file.
# (Touching an unrelated toggle property forces the host to keep the
# touched segments as distinct runs instead of re-merging them with their
# neighbours once the colour/highlight above made them all identical.)
$segments = @("https://", "www.internetingishard.com", "/html-and-", "css", "/basic-web-pages/")
$pos = $urlPara.Start
for ($i = 0; $i -lt $segments.Count; $i++) {
    $len = $segments[$i].Length
    if ($i % 2 -eq 1) {
        $seg = $tr.Characters($pos, $len)
        $seg.Font.Italic = -1
        $seg.Font.Italic = 0
    }
    $pos = $pos + $len
}

# 2) Now give the trailing "html" run the same hyperlink as the w3schools
#    run right before it (added last so it doesn't bleed into the newly
#    inserted paragraphs above).
$htmlRun = $seePara.Runs(2)
$htmlRun.ActionSettings(1).Hyperlink.Address = "https://www.w3schools.com/"
